$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column E (the 5th column). This shifts the
# existing "echo" / "executeAsyncScript" / "wait" test-case columns one
# slot to the right (E->F, F->G, G->H) and copies column D's formatting
# into the freshly inserted column E.
$ws.Columns.Item(5).Insert()

# Populate the new column E with a second "executeScript" test case that
# uses the JSON-wrapped parameter/result format.
$ws.Range("E1").Value = "executeScript"
$ws.Range("E2").Value = '{"target":"return document.title;"}'
$ws.Range("E3").Value = '{"value":"result"}'

# The "${result}" placeholder used to live in column E (now shifted to F)
# in row 3. In the new layout it instead belongs in row 2 of column F, so
# move it up a row and give it the same formatting as the other row-2
# parameter cells in this column group (copied from D3, which carries
# that style).
$ws.Range("D3").Copy()
$ws.Range("F2").PasteSpecial(-4122)
$ws.Range("F2").Value = '${result}'

# Row 3 no longer has any content/formatting carried in column F.
$ws.Range("F3").Clear()

# Column E should be noticeably wider than the other parameter columns.
$ws.Columns.Item(5).ColumnWidth = 23.2857142857

$null = $ws.Range("F2").Select()
